$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values (quarter 01-01-2021)
$ws.Range("B74").Value = -1836.8
$ws.Range("C74").Value = 4265.5
$ws.Range("D74").Value = 124.7
$ws.Range("E74").Value = 11234.4

# Add new row 75 (quarter 01-04-2021)
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").ClearFormats()
$ws.Range("B75").Value = -2637.7
$ws.Range("C75").Value = 4174.2
$ws.Range("D75").Value = -12448
$ws.Range("E75").Value = 13972.8
